$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell "Save" in H1, reusing the same formatting as the other
# header cells (bold, bordered, centered) by copying G1's format.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New "Save" column values for rows 2-9
$values = @(1, 0, 1, 1, 1, 1, 0, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
